# Append new tyre sales rows (278-300) to the "Holidays 2019" sheet,
# mirroring the season-filtered Onliner price data added on 2022-12-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel date serial epoch (serial 1 = 1899-12-31, with the classic 1900 leap-year bug,
# so day 0 == 1899-12-30) used to turn the raw serials from the source data into
# real dates that Excel will store/format correctly.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

# Columns: RowNumber, E=Tyre Size, F=Model, G=Param, H=Sales value, I=Date serial, J=Contragent
$data = @(
    @(278, "35/65-33", "ФБел-283", "42, груз, сер", 2, 44910, "нет данных"),
    @(279, "205/55R16", "BEL-262", "б/к, сер, легк", 2, 44910, "нет данных"),
    @(280, "205/55R16", "BEL-317", "б/к, сер, легк", 2, 44910, "нет данных"),
    @(281, "205/55R16", "BEL-317S", "сер, ошип", 2, 44910, "нет данных"),
    @(282, "24.00R35", "Бел-122", "груз, сер, LS-2, Type", 2, 44910, "нет данных"),
    @(283, "24.00R35", "Бел-122", "груз, сер, Type, H", 2, 44910, "нет данных"),
    @(284, "24.00R35", "Бел-122", "груз, сер, Type, C", 2, 44910, "нет данных"),
    @(285, "24.00R35", "Бел-202", "210B, сер, LS-2, Type", 2, 44910, "нет данных"),
    @(286, "24.00R35", "Бел-202", "210B, сер, Type, H", 2, 44910, "нет данных"),
    @(287, "24.00R35", "Бел-202", "210B, сер, Type, C", 2, 44910, "нет данных"),
    @(288, "24.00R35", "Бел-212", "груз, сер, LS-2, Type", 2, 44910, "нет данных"),
    @(289, "24.00R35", "Бел-212", "груз, сер, Type, H", 2, 44910, "нет данных"),
    @(290, "24.00R35", "Бел-212", "груз, сер, Type, C", 2, 44910, "нет данных"),
    @(291, "21.00R35", "Бел-200", "202B, сер, LS-2, Type", 2, 44910, "нет данных"),
    @(292, "21.00R35", "Бел-200", "202B, сер, Type, H", 2, 44910, "нет данных"),
    @(293, "21.00R35", "Бел-200", "202B, сер, Type, C", 2, 44910, "нет данных"),
    @(294, "21.00R35", "Бел-210", "202B, сер, LS-2, Type", 2, 44910, "нет данных"),
    @(295, "21.00R35", "Бел-210", "202B, сер, Type, H", 2, 44910, "нет данных"),
    @(296, "21.00R35", "Бел-210", "202B, сер, Type, C", 2, 44910, "нет данных"),
    @(297, "14.00R20", "BEL-248", "груз, сер", 2, 44910, "нет данных"),
    @(298, "14.00R20", "BEL-248", "груз, сер", 2, 44910, "нет данных"),
    @(299, "14.00R20", "BEL-248", "б/к, груз, сер", 2, 44910, "нет данных"),
    @(300, "14.00R20", "BEL-248", "груз, сер", 2, 44910, "нет данных")
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
    $ws.Cells.Item($r, 8).Value = $row[4]

    # Date column: set the display format first so the stored style matches
    # the yyyy-mm-dd format already used by the rest of column I.
    $ws.Cells.Item($r, 9).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 9).Value = $epoch.AddDays($row[5])

    $ws.Cells.Item($r, 10).Value = $row[6]
}
